$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# Final target data for rows 2-15 (columns A-H):
# A=Client Id, B=Candidate ID, C=User Name, D=Exam Password,
# E=Title, F=First Name, G=Last Name, H=Role
$data = @(
    ,(@('qfwoq293', 231011212, 'dryduis30', 'xM8d4A$#', 'MR', 'sHxFDoDV', 'XNjp', 'Candidate'))
    ,(@('Gnzrc374', 231011211, 'ghhwpjt89', 'R4!E&ge6', 'MR', 'hTAAMPQQ', 'Wygx', 'Candidate'))
    ,(@('oKQVi213', 231011210, 'qxzmswk97', 'z#7qK5%Z', 'MR', 'qSxbOOcp', 'XJDD', 'Candidate'))
    ,(@('StiUU984', 231011209, 'zmdjric82', 'jqB%T79!', 'MR', 'ydFHZsTK', 'NUIm', 'Candidate'))
    ,(@('zezyW344', 231011208, 'sjncgxx74', 'u!6#p8CK', 'MR', 'xdhzjxyb', 'qBTU', 'Candidate'))
    ,(@('zQKYj409', 231011207, 'exgkopa11', 'pv5#%Z7E', 'MR', 'RZRkSMCD', 'ZEiB', 'Candidate'))
    ,(@('sdvdz650', 231011206, 'qbljmxe43', 'Qk%6&e2E', 'MR', 'DGVBCXrr', 'wLFL', 'Candidate'))
    ,(@('VPuqO795', 231011205, 'ugoeyxi89', 'F%9!V7bu', 'MR', 'WpeBGHGN', 'XueW', 'Candidate'))
    ,(@('lFyII457', 231011204, 'pmndxre61', 'Mh3Xq%7$', 'MR', 'JzGREcIy', 'RbaQ', 'Candidate'))
    ,(@('BQvSQ295', 231011203, 'vtmvvvk82', 'B6nX!3%u', 'MR', 'xfPmOWFo', 'DaKC', 'Candidate'))
    ,(@('cqTCl395', 231011202, 'qunwrpt14', 'wV6!t$W4', 'MR', 'AJQkLcVZ', 'dtJa', 'Candidate'))
    ,(@('ifjqc544', 231011201, 'ccnqbky58', 'arV#2E%8', 'MR', 'AtxritWN', 'czkw', 'Candidate'))
    ,(@('VhclZ700', 231011200, 'mprggfi15', 'a9%Vh2E&', 'MR', 'PvIVneGS', 'YcRZ', 'Candidate'))
    ,(@('Zxdni552', 231011199, 'zfnakcf63', 'y!c6%9RP', 'MR', 'gLWdkwgL', 'DxBa', 'Candidate'))
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]

    # New rows (13-15) need the same style as the existing body rows.
    if ($rowIndex -gt 12) {
        $srcRange = $ws.Range("A12:H12")
        $dstRange = $ws.Range("A$rowIndex`:H$rowIndex")
        $srcRange.Copy($dstRange)
        $ws.Cells.Item($rowIndex, 1).Value = $row[0]
        $ws.Cells.Item($rowIndex, 2).Value = $row[1]
        $ws.Cells.Item($rowIndex, 3).Value = $row[2]
        $ws.Cells.Item($rowIndex, 4).Value = $row[3]
        $ws.Cells.Item($rowIndex, 5).Value = $row[4]
        $ws.Cells.Item($rowIndex, 6).Value = $row[5]
        $ws.Cells.Item($rowIndex, 7).Value = $row[6]
        $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    }

    $rowIndex++
}

$ws.Range("A1:H15").Select()
